$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126, shifting existing rows 126:167 down to 127:168
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new data point
$ws.Cells.Item(126, 1).Value = 3
$ws.Cells.Item(126, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(126, 3).Value = "Coquimbo"
$ws.Cells.Item(126, 4).Value = 44559
$ws.Cells.Item(126, 5).Value = 5
$ws.Cells.Item(126, 6).Value = 100112010
$ws.Cells.Item(126, 7).Value = "Achicoria"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 95
$ws.Cells.Item(126, 11).Value = 6500
$ws.Cells.Item(126, 12).Value = 7000
$ws.Cells.Item(126, 13).Value = 6763
$ws.Cells.Item(126, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(126, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(126, 16).Value = 423
$ws.Cells.Item(126, 17).Value = 16
$ws.Cells.Item(126, 18).Value = "Hortaliza"

# Match the date style used by column D in the other rows
$ws.Cells.Item(126, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
